# Weekly fruit/vegetable price update: a new weekly record for
# "Feria Lagunitas de Puerto Montt - Albahaca" is inserted above the
# existing row 48, pushing all subsequent records down by one row
# (old row 48 -> new row 49, ..., old row 162 -> new row 163).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 48 (shifts row 48..162 down to 49..163)
$ws.Rows.Item(48).Insert()

# Populate the new row 48 with the new weekly data point
$ws.Cells.Item(48, 1).Value = 4
$ws.Cells.Item(48, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(48, 3).Value = "Los Lagos"
$ws.Cells.Item(48, 4).Value = 44952
$ws.Cells.Item(48, 5).Value = 10
$ws.Cells.Item(48, 6).Value = 100112052
$ws.Cells.Item(48, 7).Value = "Albahaca"
$ws.Cells.Item(48, 8).Value = "Sin especificar"
$ws.Cells.Item(48, 9).Value = "Primera"
$ws.Cells.Item(48, 10).Value = 60
$ws.Cells.Item(48, 11).Value = 6000
$ws.Cells.Item(48, 12).Value = 6000
$ws.Cells.Item(48, 13).Value = 6000
$ws.Cells.Item(48, 14).Value = "`$/docena de matas"
$ws.Cells.Item(48, 15).Value = "Región Metropolitana"
$ws.Cells.Item(48, 16).Value = 1000
$ws.Cells.Item(48, 17).Value = 6
$ws.Cells.Item(48, 18).Value = "Hortaliza"
